$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 786.8461
$ws.Range("I43").Value = 586.3333
$ws.Range("J43").Value = 847
$ws.Range("K43").Value = 586.3333
$ws.Range("L43").Value = 847
$ws.Range("M43").Value = -517.3333
$ws.Range("N43").Value = -985

$ws.Range("H62").Value = 2936.7144
$ws.Range("I62").Value = 2470
$ws.Range("J62").Value = 3361
$ws.Range("K62").Value = 2470
$ws.Range("L62").Value = 3361
$ws.Range("M62").Value = -1846
$ws.Range("N62").Value = -4609

$ws.Range("H65").Value = 2936.7144
$ws.Range("I65").Value = 2470
$ws.Range("J65").Value = 3361
$ws.Range("K65").Value = 12350
$ws.Range("L65").Value = 16805
$ws.Range("M65").Value = -9230
$ws.Range("N65").Value = -23045

$ws.Range("H92").Value = 418.3889
$ws.Range("I92").Value = 401.9375
$ws.Range("K92").Value = 401.9375
$ws.Range("M92").Value = 846.0625

$ws.Range("H103").Value = 333333340
$ws.Range("I103").Value = 333333340
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 1000000020
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -999999434
$ws.Range("N103").ClearContents()

$ws.Range("H138").Value = 2275.524
$ws.Range("I138").Value = 1203.9131
$ws.Range("J138").Value = 3572.7368
$ws.Range("K138").Value = 3611.7393
$ws.Range("L138").Value = 10718.2104
$ws.Range("M138").Value = 1528.2607
$ws.Range("N138").Value = -20998.2104

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1437.1333
$ws.Range("I102").Value = 1380.4615
$ws.Range("K102").Value = 1380.4615
$ws.Range("M102").Value = 241.5385000000001

$ws.Range("H107").Value = 19728.5
$ws.Range("J107").Value = 19728.5
$ws.Range("L107").Value = 19728.5
$ws.Range("N107").Value = -27408.5

$ws.Range("H114").Value = 31556
$ws.Range("J114").Value = 31556
$ws.Range("L114").Value = 31556
$ws.Range("N114").Value = -40234

$ws.Range("H132").Value = 13923.69
$ws.Range("I132").Value = 1770.0667
$ws.Range("J132").Value = 44307.75
$ws.Range("K132").Value = 5310.2001
$ws.Range("L132").Value = 132923.25
$ws.Range("M132").Value = -2780.2001
$ws.Range("N132").Value = -137983.25

$ws.Range("H139").Value = 40423.418
$ws.Range("J139").Value = 40423.418
$ws.Range("L139").Value = 40423.418
$ws.Range("N139").Value = -50703.418

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 50000256
$ws.Range("I64").Value = 125000130
$ws.Range("J64").Value = 342.66666
$ws.Range("K64").Value = 125000130
$ws.Range("L64").Value = 342.66666
$ws.Range("M64").Value = -124999905
$ws.Range("N64").Value = -792.66666

$ws.Range("H67").Value = 50000256
$ws.Range("I67").Value = 125000130
$ws.Range("J67").Value = 342.66666
$ws.Range("K67").Value = 125000130
$ws.Range("L67").Value = 342.66666
$ws.Range("M67").Value = -124999350
$ws.Range("N67").Value = -1902.66666

$ws.Range("H105").Value = 1564315.2
$ws.Range("I105").Value = 1660.5385
$ws.Range("K105").Value = 1660.5385
$ws.Range("M105").Value = 86.46149999999989

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 16583.758
$ws.Range("I58").Value = 1480.9445
$ws.Range("J58").Value = 34707.133
$ws.Range("K58").Value = 1480.9445
$ws.Range("L58").Value = 34707.133
$ws.Range("M58").Value = -1277.9445
$ws.Range("N58").Value = -35113.133

$ws.Range("H120").Value = 13666.667
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 13666.667
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 13666.667
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -20924.667

$ws.Range("H121").Value = 9604
$ws.Range("I121").Value = 5425
$ws.Range("J121").Value = 16290.4
$ws.Range("K121").Value = 5425
$ws.Range("L121").Value = 16290.4
$ws.Range("M121").Value = -4115
$ws.Range("N121").Value = -18910.4

$ws.Range("H134").Value = 1368.8235
$ws.Range("I134").Value = 1174.6154
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 3523.8462
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -988.8462
$ws.Range("N134").Value = -11070

$ws.Range("H136").Value = 16583.758
$ws.Range("I136").Value = 1480.9445
$ws.Range("J136").Value = 34707.133
$ws.Range("K136").Value = 4442.833500000001
$ws.Range("L136").Value = 104121.399
$ws.Range("M136").Value = -1892.833500000001
$ws.Range("N136").Value = -109221.399

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H117").Value = 999.2
$ws.Range("I117").Value = 751.1667
$ws.Range("J117").Value = 1164.5555
$ws.Range("K117").Value = 2253.5001
$ws.Range("L117").Value = 3493.6665
$ws.Range("M117").Value = 1188.4999
$ws.Range("N117").Value = -10377.6665

$ws.Range("H118").Value = 100001210
$ws.Range("I118").Value = 125000260
$ws.Range("J118").Value = 5000
$ws.Range("K118").Value = 375000780
$ws.Range("L118").Value = 15000
$ws.Range("M118").Value = -374999537
$ws.Range("N118").Value = -17486

$ws.Range("H119").Value = 5142
$ws.Range("I119").Value = 5142
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 15426
$ws.Range("L119").Value = 0
$ws.Range("M119").ClearContents()
$ws.Range("N119").Value = -10588

$ws.Range("H120").Value = 13102.7
$ws.Range("I120").Value = 6205.4
$ws.Range("K120").Value = 18616.2
$ws.Range("M120").Value = -13778.2

$ws.Range("H121").Value = 7247413
$ws.Range("I121").Value = 265
$ws.Range("J121").Value = 7937617.5
$ws.Range("K121").Value = 795
$ws.Range("L121").Value = 23812852.5
$ws.Range("M121").Value = 515
$ws.Range("N121").Value = -23815472.5

$ws.Range("H131").Value = 714.39
$ws.Range("J131").Value = 741.51086
$ws.Range("L131").Value = 2224.53258
$ws.Range("N131").Value = -12304.53258

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 25326
$ws.Range("J45").Value = 25326
$ws.Range("L45").Value = 25326
$ws.Range("N45").Value = -26444

$ws.Range("H97").Value = 1395.5264
$ws.Range("I97").Value = 1342.0588
$ws.Range("J97").Value = 1850
$ws.Range("K97").Value = 1342.0588
$ws.Range("L97").Value = 1850
$ws.Range("M97").Value = -846.0588
$ws.Range("N97").Value = -2842

$ws.Range("H102").Value = 1606.3214
$ws.Range("I102").Value = 1283.96
$ws.Range("K102").Value = 1283.96
$ws.Range("M102").Value = 338.04

$ws.Range("H132").Value = 27539.545
$ws.Range("I132").Value = 5051.5884
$ws.Range("K132").Value = 15154.7652
$ws.Range("M132").Value = -12624.7652

$ws.Range("H139").Value = 27216
$ws.Range("J139").Value = 27216
$ws.Range("L139").Value = 27216
$ws.Range("N139").Value = -37496

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 2977.5
$ws.Range("J24").Value = 2977.5
$ws.Range("L24").Value = 2977.5
$ws.Range("N24").Value = -3663.5

$ws.Range("H46").Value = 1107.28
$ws.Range("I46").Value = 1082.6459
$ws.Range("J46").Value = 1698.5
$ws.Range("K46").Value = 1082.6459
$ws.Range("L46").Value = 1698.5
$ws.Range("M46").Value = -894.6459
$ws.Range("N46").Value = -2074.5

$ws.Range("H48").Value = 15520.5
$ws.Range("I48").Value = 13041
$ws.Range("K48").Value = 13041
$ws.Range("M48").Value = -12380

$ws.Range("H106").Value = 21370
$ws.Range("J106").Value = 21370
$ws.Range("L106").Value = 21370
$ws.Range("N106").Value = -23894

$ws.Range("H136").Value = 1640.5927
$ws.Range("I136").Value = 1472.9231
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 4418.7693
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -1868.7693
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 18500
$ws.Range("J86").Value = 18500
$ws.Range("L86").Value = 18500
$ws.Range("N86").Value = -20746

$ws.Range("H89").Value = 18500
$ws.Range("J89").Value = 18500
$ws.Range("L89").Value = 92500
$ws.Range("N89").Value = -103732

$ws.Range("H136").Value = 17243352
$ws.Range("I136").Value = 20000782
$ws.Range("K136").Value = 60002346
$ws.Range("M136").Value = -59999796
